$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 238-272 (values shift down one period; see commit diff) ---
# Row 238
$ws.Range("D238").Value = 44476
$ws.Range("K238").Value = 'Packham''s Triumph'
$ws.Range("M238").Value = 200
$ws.Range("O238").Value = 10000
$ws.Range("P238").Value = 9500
$ws.Range("S238").Value = 594

# Row 239
$ws.Range("D239").Value = 44476
$ws.Range("K239").Value = 'Packham''s Triumph'
$ws.Range("M239").Value = 100
$ws.Range("N239").Value = 8000
$ws.Range("O239").Value = 8000
$ws.Range("P239").Value = 8000
$ws.Range("S239").Value = 500

# Row 240
$ws.Range("K240").Value = 'Abate Fettel'

# Row 241
$ws.Range("K241").Value = 'Abate Fettel'

# Row 242
$ws.Range("D242").Value = 44384
$ws.Range("M242").Value = 50
$ws.Range("N242").Value = 9000
$ws.Range("P242").Value = 9000
$ws.Range("S242").Value = 562

# Row 243
$ws.Range("D243").Value = 44384

# Row 244
$ws.Range("K244").Value = 'Packham''s Triumph'
$ws.Range("M244").Value = 100
$ws.Range("N244").Value = 8000
$ws.Range("O244").Value = 9000
$ws.Range("P244").Value = 8500
$ws.Range("S244").Value = 531

# Row 245
$ws.Range("K245").Value = 'Packham''s Triumph'
$ws.Range("M245").Value = 50
$ws.Range("N245").Value = 7000
$ws.Range("O245").Value = 7000
$ws.Range("P245").Value = 7000
$ws.Range("S245").Value = 438

# Row 246
$ws.Range("D246").Value = 44363
$ws.Range("K246").Value = 'Winter Nelis'

# Row 247
$ws.Range("D247").Value = 44363
$ws.Range("K247").Value = 'Winter Nelis'

# Row 248
$ws.Range("D248").Value = 44306
$ws.Range("K248").Value = 'Packham''s Triumph'
$ws.Range("M248").Value = 200
$ws.Range("O248").Value = 10000
$ws.Range("P248").Value = 9500
$ws.Range("S248").Value = 594

# Row 249
$ws.Range("D249").Value = 44306
$ws.Range("K249").Value = 'Packham''s Triumph'
$ws.Range("N249").Value = 8000
$ws.Range("O249").Value = 8000
$ws.Range("P249").Value = 8000
$ws.Range("S249").Value = 500

# Row 250
$ws.Range("K250").Value = 'Abate Fettel'

# Row 251
$ws.Range("K251").Value = 'Abate Fettel'

# Row 252
$ws.Range("D252").Value = 44357
$ws.Range("K252").Value = 'Packham''s Triumph'
$ws.Range("O252").Value = 9000
$ws.Range("P252").Value = 9000
$ws.Range("S252").Value = 562

# Row 253
$ws.Range("D253").Value = 44357
$ws.Range("K253").Value = 'Packham''s Triumph'
$ws.Range("M253").Value = 100
$ws.Range("N253").Value = 7000
$ws.Range("O253").Value = 7000
$ws.Range("P253").Value = 7000
$ws.Range("S253").Value = 438

# Row 254
$ws.Range("K254").Value = 'Beurre Bosc'
$ws.Range("M254").Value = 100
$ws.Range("O254").Value = 10000
$ws.Range("P254").Value = 9500
$ws.Range("S254").Value = 594

# Row 255
$ws.Range("K255").Value = 'Beurre Bosc'

# Row 256
$ws.Range("K256").Value = 'Forelle'
$ws.Range("L256").Value = 'Primera'
$ws.Range("N256").Value = 9000
$ws.Range("O256").Value = 9000
$ws.Range("P256").Value = 9000
$ws.Range("S256").Value = 562

# Row 257
$ws.Range("K257").Value = 'Forelle'
$ws.Range("L257").Value = 'Segunda'
$ws.Range("M257").Value = 50
$ws.Range("N257").Value = 8000
$ws.Range("O257").Value = 8000
$ws.Range("P257").Value = 8000
$ws.Range("S257").Value = 500

# Row 258
$ws.Range("L258").Value = 'Especial'
$ws.Range("M258").Value = 50
$ws.Range("N258").Value = 12000
$ws.Range("O258").Value = 12000
$ws.Range("P258").Value = 12000
$ws.Range("S258").Value = 750

# Row 259
$ws.Range("D259").Value = 44328
$ws.Range("K259").Value = 'Packham''s Triumph'
$ws.Range("N259").Value = 10000
$ws.Range("P259").Value = 10000
$ws.Range("S259").Value = 625

# Row 260
$ws.Range("D260").Value = 44328
$ws.Range("K260").Value = 'Packham''s Triumph'
$ws.Range("M260").Value = 100

# Row 261
$ws.Range("D261").Value = 44321
$ws.Range("K261").Value = 'Winter Nelis'
$ws.Range("N261").Value = 9000
$ws.Range("O261").Value = 10000
$ws.Range("P261").Value = 9500
$ws.Range("S261").Value = 594

# Row 262
$ws.Range("D262").Value = 44321
$ws.Range("K262").Value = 'Winter Nelis'
$ws.Range("N262").Value = 8000
$ws.Range("O262").Value = 8000
$ws.Range("P262").Value = 8000
$ws.Range("S262").Value = 500

# Row 263
$ws.Range("D263").Value = 44397
$ws.Range("N263").Value = 8000
$ws.Range("O263").Value = 9000
$ws.Range("P263").Value = 8500
$ws.Range("S263").Value = 531

# Row 264
$ws.Range("D264").Value = 44397
$ws.Range("N264").Value = 7000
$ws.Range("O264").Value = 7000
$ws.Range("P264").Value = 7000
$ws.Range("S264").Value = 438

# Row 265
$ws.Range("D265").Value = 44314
$ws.Range("N265").Value = 9000
$ws.Range("O265").Value = 10000
$ws.Range("P265").Value = 9500
$ws.Range("S265").Value = 594

# Row 266
$ws.Range("D266").Value = 44314
$ws.Range("N266").Value = 8000
$ws.Range("O266").Value = 8000
$ws.Range("P266").Value = 8000
$ws.Range("S266").Value = 500

# Row 267
$ws.Range("D267").Value = 44392
$ws.Range("M267").Value = 100
$ws.Range("N267").Value = 8000
$ws.Range("O267").Value = 9000
$ws.Range("P267").Value = 8500
$ws.Range("S267").Value = 531

# Row 268
$ws.Range("D268").Value = 44392
$ws.Range("M268").Value = 50
$ws.Range("N268").Value = 7000
$ws.Range("O268").Value = 7000
$ws.Range("P268").Value = 7000
$ws.Range("S268").Value = 438

# Row 269
$ws.Range("D269").Value = 44425
$ws.Range("N269").Value = 9000
$ws.Range("O269").Value = 10000
$ws.Range("P269").Value = 9500
$ws.Range("S269").Value = 594

# Row 270
$ws.Range("D270").Value = 44425
$ws.Range("N270").Value = 8000
$ws.Range("O270").Value = 8000
$ws.Range("P270").Value = 8000
$ws.Range("S270").Value = 500

# Row 271
$ws.Range("D271").Value = 44390
$ws.Range("N271").Value = 8000
$ws.Range("O271").Value = 9000
$ws.Range("P271").Value = 8500
$ws.Range("S271").Value = 531

# Row 272
$ws.Range("D272").Value = 44390
$ws.Range("N272").Value = 7000
$ws.Range("O272").Value = 7000
$ws.Range("P272").Value = 7000
$ws.Range("S272").Value = 438

# --- Append new rows 273-274 (new reporting week) ---
# Row 273
$ws.Range("A273").Value = 11
$ws.Range("B273").Value = 'Vega Monumental Concepción'
$ws.Range("C273").Value = 'Bíobío'
$ws.Range("D273").Value = 44432
$ws.Range("E273").Value = 8
$ws.Range("F273").Value = 'Fruta'
$ws.Range("G273").Value = 100104
$ws.Range("H273").Value = 'Frutos de pepita'
$ws.Range("I273").Value = 100104005
$ws.Range("J273").Value = 'Pera'
$ws.Range("K273").Value = 'Packham''s Triumph'
$ws.Range("L273").Value = 'Primera'
$ws.Range("M273").Value = 200
$ws.Range("N273").Value = 9000
$ws.Range("O273").Value = 10000
$ws.Range("P273").Value = 9500
$ws.Range("Q273").Value = '$/caja 16 kilos empedrada'
$ws.Range("R273").Value = 'Región de O''Higgins'
$ws.Range("S273").Value = 594
$ws.Range("T273").Value = 16
$ws.Range("D273").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 274
$ws.Range("A274").Value = 11
$ws.Range("B274").Value = 'Vega Monumental Concepción'
$ws.Range("C274").Value = 'Bíobío'
$ws.Range("D274").Value = 44432
$ws.Range("E274").Value = 8
$ws.Range("F274").Value = 'Fruta'
$ws.Range("G274").Value = 100104
$ws.Range("H274").Value = 'Frutos de pepita'
$ws.Range("I274").Value = 100104005
$ws.Range("J274").Value = 'Pera'
$ws.Range("K274").Value = 'Packham''s Triumph'
$ws.Range("L274").Value = 'Segunda'
$ws.Range("M274").Value = 100
$ws.Range("N274").Value = 8000
$ws.Range("O274").Value = 8000
$ws.Range("P274").Value = 8000
$ws.Range("Q274").Value = '$/caja 16 kilos empedrada'
$ws.Range("R274").Value = 'Región de O''Higgins'
$ws.Range("S274").Value = 500
$ws.Range("T274").Value = 16
$ws.Range("D274").NumberFormat = "YYYY-MM-DD HH:MM:SS"
